$wb = $excel.ActiveWorkbook

# --- Section_A (Sheet1) updates ---
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("C2").Value = "EC301"
$wsA.Range("D2").Value = "MA261"
$wsA.Range("E2").Value = "EC302"
$wsA.Range("F2").Value = "HS201"

$wsA.Range("C3").Value = "CS263"
$wsA.Range("D3").Value = "HS201"
$wsA.Range("E3").Value = "CS263"
$wsA.Range("F3").Value = "EC302"

$wsA.Range("B5").Value = "EC301"
$wsA.Range("C5").Value = "MA261"
$wsA.Range("D5").Value = "EC301"
$wsA.Range("E5").Value = "CS251 (Elective)"
$wsA.Range("F5").Value = "CS263"

$wsA.Range("C6").Value = "Free"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "CS263 (Tutorial)"
$wsA.Range("F6").Value = "EC301 (Tutorial)"

$wsA.Range("C7").Value = "HS201"
$wsA.Range("F7").Value = "MA262"

$wsA.Range("C8").Value = "EC302 (Tutorial)"
$wsA.Range("D8").Value = "Free"
$wsA.Range("F8").Value = "Free"

# --- Section_B (Sheet2) updates ---
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("C2").Value = "HS201"
$wsB.Range("D2").Value = "MA261"
$wsB.Range("E2").Value = "CS263"
$wsB.Range("F2").Value = "EC302"

$wsB.Range("C3").Value = "EC301"
$wsB.Range("D3").Value = "HS201"
$wsB.Range("E3").Value = "EC301"

$wsB.Range("B5").Value = "EC302"
$wsB.Range("D5").Value = "MA262"
$wsB.Range("E5").Value = "CS251 (Elective)"
$wsB.Range("F5").Value = "EC301"

$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "CS263 (Tutorial)"
$wsB.Range("E6").Value = "EC301 (Tutorial)"

$wsB.Range("B7").Value = "HS201"
$wsB.Range("C7").Value = "EC302"
$wsB.Range("D7").Value = "CS263"
$wsB.Range("F7").Value = "MA261"

$wsB.Range("B8").Value = "EC302 (Tutorial)"
$wsB.Range("C8").Value = "Free"
$wsB.Range("D8").Value = "Free"

$wb.Save()
